# by jiankong on 1113
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: drop the two empty "张悦" / "卢楠" rows, keep the real data rows ---
$ws1.Range("A2:A3").EntireRow.Delete()

# --- Rename the group from "质控组" to "北京组" on the remaining Sheet1 rows ---
$ws1.Range("A2").Value = "北京组"
$ws1.Range("A3").Value = "北京组"

# --- Sheet2: rename the group (member/aggregate label "总体" stays as-is) ---
$ws2.Range("A2").Value = "北京组"

# --- Update selections on each sheet ---
$ws1.Range("C11").Select()
$ws2.Range("A2").Select()

# --- Make Sheet2 the active tab (matches tabSelected moving to Sheet2) ---
$ws2.Activate()
